$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-metric table values (ifo GDP component analysis preprocessing).
# Rows correspond to quarter buckets Q0..Q8 (row 2..10); columns B:F are
# ME, MAE, MSE, RMSE, SE and column G is the observation count N.

$data = @{
    2  = @(0.04628788019112266, 0.4824692687926074, 0.3511399927416267, 0.5925706647663439, 0.606948243852973, 19)
    3  = @(0.4480479793330983, 0.7418921356350544, 0.8887687429572028, 0.9427453224265835, 0.8535183764212539, 18)
    4  = @(0.6602403078873996, 0.8862161360526695, 1.515028269969415, 1.230864846345615, 1.070773292379532, 17)
    5  = @(0.6733686636360297, 0.7201182210757808, 0.9403137650122958, 0.9696977699326197, 0.7206577794175344, 16)
    6  = @(0.5133116504495733, 0.6631082909423874, 0.7111133683316896, 0.84327538107767, 0.692529925483697, 15)
    7  = @(0.458951976745339, 0.673915242314678, 0.6571042978108756, 0.8106196998660196, 0.6934053832595887, 14)
    8  = @(0.3510070987759771, 0.6341193177088349, 0.5711820071072622, 0.7557658414530669, 0.696640049829043, 13)
    9  = @(0.6367582691617952, 0.6367582691617952, 0.6187545134813729, 0.7866095050794727, 0.4898521181462611, 9)
    10 = @(0.4167169463657505, 0.4167169463657505, 0.2338164441105596, 0.4835457001262234, 0.2742340030023716, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $col = 2
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}
